$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.904.99'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '2.785.90'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '356.72'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '109.55'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.24%  '
$ws.Range('E7').Value = '  -1.81%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.588'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.81%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.31'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('E11').Value = '  +1.51%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0847'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.11%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.48'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.98%  '
$ws.Range('E14').Value = '  -3.04%  '
$ws.Range('D15').Value = '3.228.18'
$ws.Range('E15').Value = '  -1.85%  '
$ws.Range('D16').Value = '2.817.15'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.946'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.52%  '
$ws.Range('D18').Value = '51.878.46'
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.48'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.08'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.54%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.15'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.59%  '
$ws.Range('D22').Value = '0.0₃0975'
$ws.Range('E22').Value = '  -2.49%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.24'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '270.41'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('E25').Value = '  -4.84%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.48'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.164'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +16.44%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.31'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.11%  '
$ws.Range('E30').Value = '  -5.09%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0471'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '52.04'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.64'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.74'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0844'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.59%  '
$ws.Range('E36').Value = '  -6.47%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.79'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.05%  '
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.98'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.59%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.60'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.82%  '
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '119.75'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.14%  '
$ws.Range('E45').Value = '  -8.06%  '
$ws.Range('D46').Value = '2.081.60'
$ws.Range('E46').Value = '  -1.54%  '
$ws.Range('E47').Value = '  -4.55%  '
$ws.Range('E48').Value = '  -2.00%  '
$ws.Range('E49').Value = '  -5.20%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.953'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.40%  '
$ws.Range('E51').Value = '  +32.57%  '
